$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Ich bin nicht sicher, ob Klimawissenschaftler ihr Wissen vollständig mit der Öffentlichkeit teilen."
$ws.Range("B3").Value = "Der beobachtete Klimawandel ist ausschließlich ein natürlicher Prozess."
$ws.Range("B4").Value = "Ich denke, der Klimawandel ist ein ernstes Problem."
$ws.Range("B5").Value = "Die Menschheit kann nicht viel tun, um die Umweltprobleme zu lösen."
$ws.Range("B6").Value = "Ich glaube, dass die meisten Umweltprobleme übertrieben werden."
$ws.Range("B7").Value = "Die Menschheit ist maßgeblich für die globale Erwärmung verantwortlich."
$ws.Range("B8").Value = "Ich glaube, dass die meisten Behauptungen über den Klimawandel wahr sind."
$ws.Range("B9").Value = "Der Versuch, Umweltprobleme zu lösen, ist reine Zeitverschwendung."
$ws.Range("B10").Value = "Ich bezweifle, dass menschliches Handeln die Erderwärmung verursacht haben."
$ws.Range("B11").Value = "Menschliches Handeln hat wenig Einfluss auf die Erderwärmung."
$ws.Range("B12").Value = "Ich bin mir nicht sicher, ob die Erderwärmung tatsächlich stattfindet."
$ws.Range("B13").Value = "Ich mache mir über die Folgen des Klimawandels Sorgen."

$ws.Range("B14").Select()
